$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 67.96296
$ws.Range("I2").Value = 68.59999999999999
$ws.Range("J2").Value = 60
$ws.Range("K2").Value = 68.59999999999999
$ws.Range("L2").Value = 60
$ws.Range("M2").Value = 44.40000000000001
$ws.Range("N2").Value = -286
$ws.Range("H17").Value = 1316.2333
$ws.Range("J17").Value = 1404.7307
$ws.Range("L17").Value = 4214.1921
$ws.Range("N17").Value = -4550.1921
$ws.Range("H32").Value = 1138
$ws.Range("I32").Value = 666
$ws.Range("J32").Value = 1798.8
$ws.Range("K32").Value = 666
$ws.Range("L32").Value = 1798.8
$ws.Range("M32").Value = -340
$ws.Range("N32").Value = -2450.8
$ws.Range("H40").Value = 3336993.2
$ws.Range("I40").Value = 3033.3333
$ws.Range("J40").Value = 4765833
$ws.Range("K40").Value = 3033.3333
$ws.Range("L40").Value = 4765833
$ws.Range("M40").Value = -2858.3333
$ws.Range("N40").Value = -4766183
$ws.Range("H97").Value = 20011600
$ws.Range("J97").Value = 20011600
$ws.Range("L97").Value = 60034800
$ws.Range("N97").Value = -60035792
$ws.Range("H98").Value = 3307
$ws.Range("I98").Value = 3307
$ws.Range("K98").Value = 3307
$ws.Range("M98").Value = -1809
$ws.Range("I107").Value = 6945319
$ws.Range("K107").Value = 6945319
$ws.Range("M107").Value = -6943399
$ws.Range("H112").Value = 4569.7144
$ws.Range("I112").Value = 1542.5
$ws.Range("J112").Value = 4802.577
$ws.Range("K112").Value = 4627.5
$ws.Range("L112").Value = 14407.731
$ws.Range("M112").Value = -3519.5
$ws.Range("N112").Value = -16623.731
$ws.Range("H122").Value = 3307
$ws.Range("I122").Value = 3307
$ws.Range("K122").Value = 9921
$ws.Range("M122").Value = -7471
$ws.Range("H129").Value = 2050
$ws.Range("J129").Value = 2050
$ws.Range("L129").Value = 6150
$ws.Range("N129").Value = -16150
$ws.Range("H137").Value = 3668.5305
$ws.Range("I137").Value = 3620.361
$ws.Range("K137").Value = 10861.083
$ws.Range("M137").Value = -8311.082999999999
$ws.Range("H138").Value = 4128.6743
$ws.Range("I138").Value = 2661.0527
$ws.Range("J138").Value = 4544.8657
$ws.Range("K138").Value = 7983.158100000001
$ws.Range("L138").Value = 13634.5971
$ws.Range("M138").Value = -2843.158100000001
$ws.Range("N138").Value = -23914.5971
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 50349
$ws.Range("I5").Value = 100000
$ws.Range("J5").Value = 698
$ws.Range("K5").Value = 100000
$ws.Range("L5").Value = 698
$ws.Range("M5").Value = -99888
$ws.Range("N5").Value = -922
$ws.Range("H32").Value = 2539218.2
$ws.Range("I32").Value = 2823412
$ws.Range("K32").Value = 2823412
$ws.Range("M32").Value = -2823125
$ws.Range("H45").Value = 3115.0588
$ws.Range("I45").Value = 2628.5454
$ws.Range("J45").Value = 4007
$ws.Range("K45").Value = 2628.5454
$ws.Range("L45").Value = 4007
$ws.Range("M45").Value = -2251.5454
$ws.Range("N45").Value = -4761
$ws.Range("H122").Value = 3001.1333
$ws.Range("I122").Value = 2345.65
$ws.Range("J122").Value = 4312.1
$ws.Range("K122").Value = 7036.950000000001
$ws.Range("L122").Value = 12936.3
$ws.Range("M122").Value = -4586.950000000001
$ws.Range("N122").Value = -17836.3
$ws.Range("H132").Value = 5874.161
$ws.Range("I132").Value = 2179.25
$ws.Range("K132").Value = 6537.75
$ws.Range("M132").Value = -4007.75
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 50349
$ws.Range("I4").Value = 100000
$ws.Range("J4").Value = 698
$ws.Range("K4").Value = 100000
$ws.Range("L4").Value = 698
$ws.Range("M4").Value = -99885
$ws.Range("N4").Value = -928
$ws.Range("H86").Value = 71906.07000000001
$ws.Range("I86").Value = 115646
$ws.Range("J86").Value = 6296.1665
$ws.Range("K86").Value = 115646
$ws.Range("L86").Value = 6296.1665
$ws.Range("M86").Value = -114523
$ws.Range("N86").Value = -8542.166499999999
$ws.Range("H89").Value = 71906.07000000001
$ws.Range("I89").Value = 115646
$ws.Range("J89").Value = 6296.1665
$ws.Range("K89").Value = 578230
$ws.Range("L89").Value = 31480.8325
$ws.Range("M89").Value = -572614
$ws.Range("N89").Value = -42712.8325
$ws.Range("H134").Value = 4634455
$ws.Range("J134").Value = 9873.526
$ws.Range("L134").Value = 29620.578
$ws.Range("N134").Value = -34690.578
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").ClearContents()
$ws.Range("N140").Value = 0
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 7943.086
$ws.Range("I31").Value = 2230.7222
$ws.Range("J31").Value = 13991.471
$ws.Range("K31").Value = 2230.7222
$ws.Range("L31").Value = 13991.471
$ws.Range("M31").Value = -1935.7222
$ws.Range("N31").Value = -14581.471
$ws.Range("H34").Value = 7943.086
$ws.Range("I34").Value = 2230.7222
$ws.Range("J34").Value = 13991.471
$ws.Range("K34").Value = 2230.7222
$ws.Range("L34").Value = 13991.471
$ws.Range("M34").Value = -2028.7222
$ws.Range("N34").Value = -14395.471
$ws.Range("H105").Value = 4763180
$ws.Range("I105").Value = 5953142.5
$ws.Range("K105").Value = 5953142.5
$ws.Range("M105").Value = -5951395.5
$ws.Range("H122").Value = 3341.4194
$ws.Range("I122").Value = 2589.5
$ws.Range("J122").Value = 5179.4443
$ws.Range("K122").Value = 7768.5
$ws.Range("L122").Value = 15538.3329
$ws.Range("M122").Value = -5318.5
$ws.Range("N122").Value = -20438.3329
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 96079.05
$ws.Range("I2").Value = 39.77778
$ws.Range("J2").Value = 168108.5
$ws.Range("K2").Value = 238.66668
$ws.Range("L2").Value = 1008651
$ws.Range("M2").Value = -125.66668
$ws.Range("N2").Value = -1008877
$ws.Range("H5").Value = 2141.0344
$ws.Range("J5").Value = 3360.6924
$ws.Range("L5").Value = 10082.0772
$ws.Range("N5").Value = -10306.0772
$ws.Range("H80").Value = 41670824
$ws.Range("I80").Value = 27781612
$ws.Range("J80").Value = 83338460
$ws.Range("K80").Value = 83344836
$ws.Range("L80").Value = 250015380
$ws.Range("M80").Value = -83343900
$ws.Range("N80").Value = -250017252
$ws.Range("H83").Value = 41670824
$ws.Range("I83").Value = 27781612
$ws.Range("J83").Value = 83338460
$ws.Range("K83").Value = 250034508
$ws.Range("L83").Value = 750046140
$ws.Range("M83").Value = -250029828
$ws.Range("N83").Value = -750055500
$ws.Range("H122").Value = 1089330.9
$ws.Range("I122").Value = 1489810.5
$ws.Range("K122").Value = 13408294.5
$ws.Range("M122").Value = -13405844.5
$ws.Range("H135").Value = 2141.0344
$ws.Range("J135").Value = 3360.6924
$ws.Range("L135").Value = 30246.2316
$ws.Range("N135").Value = -35316.2316
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 225.66667
$ws.Range("I2").Value = 39.333332
$ws.Range("K2").Value = 39.333332
$ws.Range("M2").Value = 73.666668
$ws.Range("H70").Value = 18010.87
$ws.Range("I70").Value = 21650.133
$ws.Range("K70").Value = 21650.133
$ws.Range("M70").Value = -21380.133
$ws.Range("H73").Value = 18010.87
$ws.Range("I73").Value = 21650.133
$ws.Range("K73").Value = 21650.133
$ws.Range("M73").Value = -20714.133
$ws.Range("H102").Value = 1918.6666
$ws.Range("I102").Value = 1912.8572
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1912.8572
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = -290.8571999999999
$ws.Range("N102").Value = -5244
$ws.Range("H122").Value = 4263946.5
$ws.Range("I122").Value = 4831543.5
$ws.Range("K122").Value = 14494630.5
$ws.Range("M122").Value = -14492180.5
$ws.Range("H126").Value = 7884.4653
$ws.Range("I126").Value = 5202.9375
$ws.Range("K126").Value = 15608.8125
$ws.Range("M126").Value = -13138.8125
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5832.143
$ws.Range("I7").Value = 5205
$ws.Range("J7").Value = 7400
$ws.Range("K7").Value = 5205
$ws.Range("L7").Value = 7400
$ws.Range("M7").Value = -5093
$ws.Range("N7").Value = -7624
$ws.Range("H22").Value = 6496160
$ws.Range("J22").Value = 8931957
$ws.Range("L22").Value = 8931957
$ws.Range("N22").Value = -8932547
$ws.Range("H27").Value = 6496160
$ws.Range("J27").Value = 8931957
$ws.Range("L27").Value = 8931957
$ws.Range("N27").Value = -8932171
$ws.Range("H40").Value = 4776.879
$ws.Range("I40").Value = 4208.6294
$ws.Range("J40").Value = 7334
$ws.Range("K40").Value = 4208.6294
$ws.Range("L40").Value = 7334
$ws.Range("M40").Value = -4072.6294
$ws.Range("N40").Value = -7606
$ws.Range("H46").Value = 2048.75
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 2048.75
$ws.Range("K46").Value = 0
$ws.Range("L46").ClearContents()
$ws.Range("M46").Value = 2048.75
$ws.Range("N46").Value = -2424.75
$ws.Range("H122").Value = 4797.5
$ws.Range("I122").Value = 3731.818
$ws.Range("K122").Value = 11195.454
$ws.Range("M122").Value = -8745.454000000002
$ws.Range("H126").Value = 5832.143
$ws.Range("I126").Value = 5205
$ws.Range("J126").Value = 7400
$ws.Range("K126").Value = 15615
$ws.Range("L126").Value = 22200
$ws.Range("M126").Value = -13145
$ws.Range("N126").Value = -27140
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H55").Value = 15042.286
$ws.Range("J55").Value = 22749.75
$ws.Range("L55").Value = 22749.75
$ws.Range("N55").Value = -23303.75
$ws.Range("H122").Value = 114236.53
$ws.Range("I122").Value = 183453.5
$ws.Range("K122").Value = 550360.5
$ws.Range("M122").Value = -547910.5
$ws.Range("H126").Value = 500
$ws.Range("I126").Value = 500
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 1500
$ws.Range("L126").ClearContents()
$ws.Range("N126").Value = 0
$ws.Range("M126").Value = 970
